$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset was missing a row for "L'implosione del sottomarino Titan" / FanPage / Instagram.
# Insert a new row at position 13 to restore it, which shifts all subsequent rows down by one
# (adding a new row 38 at the end), and then re-normalize every num_commenti value to 100.
$ws.Rows(13).Insert()

$titles = @(
    "Incidente Youtubers",
    "L'implosione del sottomarino Titan",
    "L'omicidio di Giulia Cecchettin",
    "Strage di Cutro"
)
$papers = @("FanPage", "Il Corriere Della Sera", "La Repubblica")
$socials = @("Facebook", "Instagram", "YouTube")

$row = 3
foreach ($titolo in $titles) {
    foreach ($giornale in $papers) {
        foreach ($social in $socials) {
            $ws.Cells.Item($row, 2).Value = $titolo
            $ws.Cells.Item($row, 3).Value = $giornale
            $ws.Cells.Item($row, 4).Value = $social
            $ws.Cells.Item($row, 5).Value = 100
            $row = $row + 1
        }
    }
}
